# Auto-generated edit script applying the Coeurl_Profits.xlsx diff
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 93.041664
$ws.Range("I92").Value = 63.88889
$ws.Range("J92").Value = 180.5
$ws.Range("K92").Value = 63.88889
$ws.Range("L92").Value = 180.5
$ws.Range("M92").Value = 1184.11111
$ws.Range("N92").Value = -2676.5
$ws.Range("H96").Value = 305.33334
$ws.Range("I96").Value = 326.75
$ws.Range("J96").Value = 262.5
$ws.Range("K96").Value = 980.25
$ws.Range("L96").Value = 787.5
$ws.Range("M96").Value = 392.75
$ws.Range("N96").Value = -3533.5
$ws.Range("H99").Value = 6702
$ws.Range("I99").Value = 2963.8
$ws.Range("J99").Value = 11374.75
$ws.Range("K99").Value = 8891.400000000001
$ws.Range("L99").Value = 34124.25
$ws.Range("M99").Value = -7393.400000000001
$ws.Range("N99").Value = -37120.25
$ws.Range("H100").Value = 2452.4783
$ws.Range("I100").Value = 2284.9167
$ws.Range("J100").Value = 2635.2727
$ws.Range("K100").Value = 2284.9167
$ws.Range("L100").Value = 2635.2727
$ws.Range("M100").Value = -1743.9167
$ws.Range("N100").Value = -3717.2727
$ws.Range("H101").Value = 323.57144
$ws.Range("J101").Value = 1000
$ws.Range("L101").Value = 3000
$ws.Range("N101").Value = -6244
$ws.Range("H129").Value = 3835.45
$ws.Range("I129").Value = 2079.6875
$ws.Range("K129").Value = 6239.0625
$ws.Range("M129").Value = -1239.0625
$ws.Range("H137").Value = 10677.786
$ws.Range("I137").Value = 2531.6667
$ws.Range("J137").Value = 16787.375
$ws.Range("K137").Value = 7595.000100000001
$ws.Range("L137").Value = 50362.125
$ws.Range("M137").Value = -5045.000100000001
$ws.Range("N137").Value = -55462.125
$ws.Range("H138").Value = 4333.829
$ws.Range("I138").Value = 1910
$ws.Range("K138").Value = 5730
$ws.Range("M138").Value = -590

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 29500
$ws.Range("J28").Value = 50000
$ws.Range("L28").Value = 50000
$ws.Range("N28").Value = -50384
$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("L46").ClearContents()
$ws.Range("N46").Value = 0
$ws.Range("H97").Value = 847.7273
$ws.Range("I97").Value = 866.3158
$ws.Range("J97").Value = 730
$ws.Range("K97").Value = 866.3158
$ws.Range("L97").Value = 730
$ws.Range("M97").Value = -370.3158
$ws.Range("N97").Value = -1722
$ws.Range("H99").Value = 29500
$ws.Range("J99").Value = 50000
$ws.Range("L99").Value = 50000
$ws.Range("N99").Value = -55990
$ws.Range("H104").Value = 149989
$ws.Range("J104").Value = 149989
$ws.Range("L104").Value = 149989
$ws.Range("N104").Value = -156977

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 6608.625
$ws.Range("I86").Value = 6834.6
$ws.Range("J86").Value = 6232
$ws.Range("K86").Value = 6834.6
$ws.Range("L86").Value = 6232
$ws.Range("M86").Value = -5711.6
$ws.Range("N86").Value = -8478
$ws.Range("H89").Value = 6608.625
$ws.Range("I89").Value = 6834.6
$ws.Range("J89").Value = 6232
$ws.Range("K89").Value = 34173
$ws.Range("L89").Value = 31160
$ws.Range("M89").Value = -28557
$ws.Range("N89").Value = -42392
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").ClearContents()
$ws.Range("N93").Value = 0
$ws.Range("H94").Value = 1579.2778
$ws.Range("I94").Value = 1386.0625
$ws.Range("K94").Value = 1386.0625
$ws.Range("M94").Value = -935.0625
$ws.Range("H99").Value = 3161.111
$ws.Range("I99").Value = 2755.9092
$ws.Range("J99").Value = 4944
$ws.Range("K99").Value = 2755.9092
$ws.Range("L99").Value = 4944
$ws.Range("M99").Value = -1257.9092
$ws.Range("N99").Value = -7940
$ws.Range("H134").Value = 3204.8696
$ws.Range("I134").Value = 2621.2666
$ws.Range("K134").Value = 7863.7998
$ws.Range("M134").Value = -5328.7998

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8825817
$ws.Range("J31").Value = 31149.46
$ws.Range("L31").Value = 31149.46
$ws.Range("N31").Value = -31739.46
$ws.Range("H34").Value = 8825817
$ws.Range("J34").Value = 31149.46
$ws.Range("L34").Value = 31149.46
$ws.Range("N34").Value = -31553.46
$ws.Range("H39").Value = 3874.75
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()
$ws.Range("H41").Value = 10999.5
$ws.Range("J41").Value = 10999.5
$ws.Range("L41").Value = 10999.5
$ws.Range("N41").Value = -11855.5
$ws.Range("H49").Value = 3874.75
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").ClearContents()
$ws.Range("H69").Value = 50962.5
$ws.Range("J69").Value = 73900
$ws.Range("L69").Value = 73900
$ws.Range("N69").Value = -75398
$ws.Range("H72").Value = 50962.5
$ws.Range("J72").Value = 73900
$ws.Range("L72").Value = 221700
$ws.Range("N72").Value = -229188

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 4622391
$ws.Range("I4").Value = 4963966
$ws.Range("K4").Value = 14891898
$ws.Range("M4").Value = -14891786
$ws.Range("H19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("L19").ClearContents()
$ws.Range("N19").Value = 0
$ws.Range("H32").Value = 800
$ws.Range("J32").Value = 800
$ws.Range("L32").Value = 2400
$ws.Range("N32").Value = -2966
$ws.Range("H86").Value = 889.6
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("H89").Value = 889.6
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").ClearContents()
$ws.Range("N106").Value = 0
$ws.Range("H122").Value = 2037.4
$ws.Range("I122").Value = 815.9
$ws.Range("K122").Value = 7343.099999999999
$ws.Range("M122").Value = -4893.099999999999

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H55").Value = 11676.667
$ws.Range("I55").Value = 11676.667
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 11676.667
$ws.Range("L55").Value = 0
$ws.Range("M55").ClearContents()
$ws.Range("N55").Value = -11349.667
$ws.Range("H62").Value = 46000
$ws.Range("I62").Value = 46000
$ws.Range("K62").Value = 46000
$ws.Range("M62").Value = -45314
$ws.Range("H65").Value = 46000
$ws.Range("I65").Value = 46000
$ws.Range("K65").Value = 138000
$ws.Range("M65").Value = -134568
$ws.Range("H70").Value = 9551.923000000001
$ws.Range("I70").Value = 11312.833
$ws.Range("J70").Value = 8042.5713
$ws.Range("K70").Value = 11312.833
$ws.Range("L70").Value = 8042.5713
$ws.Range("M70").Value = -11042.833
$ws.Range("N70").Value = -8582.5713
$ws.Range("H73").Value = 9551.923000000001
$ws.Range("I73").Value = 11312.833
$ws.Range("J73").Value = 8042.5713
$ws.Range("K73").Value = 11312.833
$ws.Range("L73").Value = 8042.5713
$ws.Range("M73").Value = -10376.833
$ws.Range("N73").Value = -9914.5713
$ws.Range("H80").Value = 2551.8333
$ws.Range("I80").Value = 1326.25
$ws.Range("K80").Value = 1326.25
$ws.Range("M80").Value = -328.25
$ws.Range("H83").Value = 2551.8333
$ws.Range("I83").Value = 1326.25
$ws.Range("K83").Value = 6631.25
$ws.Range("M83").Value = -1639.25

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 5508.5625
$ws.Range("I61").Value = 4779.769
$ws.Range("K61").Value = 4779.769
$ws.Range("M61").Value = -4577.769
$ws.Range("H93").Value = 429312.06
$ws.Range("I93").Value = 654963.75
$ws.Range("K93").Value = 654963.75
$ws.Range("M93").Value = -653715.75
$ws.Range("H100").Value = 47944.184
$ws.Range("I100").Value = 55078.527
$ws.Range("K100").Value = 55078.527
$ws.Range("M100").Value = -54537.527
$ws.Range("H113").Value = 5508.5625
$ws.Range("I113").Value = 4779.769
$ws.Range("K113").Value = 4779.769
$ws.Range("M113").Value = -2609.769
$ws.Range("H132").Value = 4741.609
$ws.Range("I132").Value = 3768.0588
$ws.Range("J132").Value = 7500
$ws.Range("K132").Value = 11304.1764
$ws.Range("L132").Value = 22500
$ws.Range("M132").Value = -8774.1764
$ws.Range("N132").Value = -27560

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H59").Value = 25000
$ws.Range("J59").Value = 25000
$ws.Range("L59").Value = 25000
$ws.Range("N59").Value = -26476
$ws.Range("H100").Value = 1448.5
$ws.Range("I100").Value = 1225.4286
$ws.Range("K100").Value = 2450.8572
$ws.Range("M100").Value = -1909.8572
$ws.Range("H107").Value = 888.34375
$ws.Range("I107").Value = 587.7778
$ws.Range("J107").Value = 1274.7858
$ws.Range("K107").Value = 1763.3334
$ws.Range("L107").Value = 3824.3574
$ws.Range("M107").Value = 156.6666
$ws.Range("N107").Value = -7664.357400000001
$ws.Range("H132").Value = 2856.4167
$ws.Range("I132").Value = 2662
$ws.Range("J132").Value = 4995
$ws.Range("K132").Value = 7986
$ws.Range("L132").Value = 14985
$ws.Range("M132").Value = -5456
$ws.Range("N132").Value = -20045
